$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 70; $row++) {
    $ws.Cells.Item($row, 24).Value = 1       # Column X = catalogo
    $ws.Cells.Item($row, 25).Value = "x"     # Column Y = descricao_tr
}
